$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '99.529.92'
$ws.Cells.Item(2, 5).Value = '  +1.00%  '
$ws.Cells.Item(3, 4).Value = '3.293.46'
$ws.Cells.Item(3, 5).Value = '  -1.81%  '
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$c = $ws.Cells.Item(5, 4)
$c.Value = "'253.63"
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.77%  '
$c = $ws.Cells.Item(6, 4)
$c.Value = "'622.97"
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.09%  '
$ws.Cells.Item(7, 5).Value = '  +18.53%  '
$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.400"
$c.Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +3.62%  '
$c = $ws.Cells.Item(10, 4)
$c.Value = "'0.968"
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +19.87%  '
$ws.Cells.Item(11, 4).Value = '3.289.90'
$ws.Cells.Item(11, 5).Value = '  -1.87%  '
$ws.Cells.Item(12, 5).Value = '  +0.10%  '
$c = $ws.Cells.Item(13, 4)
$c.Value = "'39.45"
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +10.03%  '
$ws.Cells.Item(14, 4).Value = '99.201.04'
$ws.Cells.Item(14, 5).Value = '  +0.97%  '
$c = $ws.Cells.Item(15, 4)
$c.Value = "'0.0000248"
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +0.69%  '
$ws.Cells.Item(16, 4).Value = '3.884.46'
$ws.Cells.Item(16, 5).Value = '  -2.22%  '
$c = $ws.Cells.Item(17, 4)
$c.Value = "'5.48"
$c.Style = 'Normal'
$ws.Cells.Item(18, 4).Value = '3.282.13'
$ws.Cells.Item(18, 5).Value = '  -2.18%  '
$ws.Cells.Item(19, 5).Value = '  -4.21%  '
$c = $ws.Cells.Item(20, 4)
$c.Value = "'15.39"
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +2.66%  '
$c = $ws.Cells.Item(21, 4)
$c.Value = "'6.35"
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +8.55%  '
$c = $ws.Cells.Item(22, 4)
$c.Value = "'488.72"
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +0.70%  '
$c = $ws.Cells.Item(23, 4)
$c.Value = "'9.31"
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +1.63%  '
$c = $ws.Cells.Item(24, 4)
$c.Value = "'0.0000201"
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -2.94%  '
$c = $ws.Cells.Item(25, 4)
$c.Value = "'5.65"
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -0.21%  '
$c = $ws.Cells.Item(26, 4)
$c.Value = "'89.12"
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +0.99%  '
$c = $ws.Cells.Item(27, 4)
$c.Value = "'0.323"
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +28.18%  '
$c = $ws.Cells.Item(28, 4)
$c.Value = "'12.01"
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +0.12%  '
$ws.Cells.Item(29, 4).Value = '3.434.22'
$ws.Cells.Item(29, 5).Value = '  -2.60%  '
$ws.Cells.Item(30, 5).Value = '  +0.03%  '
$ws.Cells.Item(31, 2).Value = 'Hedera'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(31, 4)
$c.Value = "'0.137"
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +9.34%  '
$ws.Cells.Item(32, 2).Value = 'Cronos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Cells.Item(32, 4)
$c.Value = "'0.190"
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +1.36%  '
$c = $ws.Cells.Item(33, 4)
$c.Value = "'10.37"
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +11.67%  '
$c = $ws.Cells.Item(34, 4)
$c.Value = "'0.999"
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +0.16%  '
$c = $ws.Cells.Item(35, 4)
$c.Value = "'27.93"
$c.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +2.26%  '
$ws.Cells.Item(36, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Cells.Item(36, 4)
$c.Value = "'0.475"
$c.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +6.25%  '
$ws.Cells.Item(37, 2).Value = 'Kaspa'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Cells.Item(37, 4)
$c.Value = "'0.149"
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -1.70%  '
$c = $ws.Cells.Item(38, 4)
$c.Value = "'7.22"
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -2.30%  '
$c = $ws.Cells.Item(39, 4)
$c.Value = "'1.93"
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -0.53%  '
$ws.Cells.Item(40, 5).Value = '  -0.22%  '
$c = $ws.Cells.Item(41, 4)
$c.Value = "'488.57"
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -5.93%  '
$ws.Cells.Item(42, 5).Value = '  -0.15%  '
$ws.Cells.Item(43, 5).Value = '  -2.24%  '
$ws.Cells.Item(44, 5).Value = '  -0.02%  '
$c = $ws.Cells.Item(45, 4)
$c.Value = "'0.773"
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -0.28%  '
$ws.Cells.Item(46, 5).Value = '  -4.73%  '
$c = $ws.Cells.Item(47, 4)
$c.Value = "'1.94"
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +0.99%  '
$ws.Cells.Item(48, 5).Value = '  -2.21%  '
$c = $ws.Cells.Item(49, 4)
$c.Value = "'0.847"
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +7.03%  '
$c = $ws.Cells.Item(50, 4)
$c.Value = "'7.28"
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +15.28%  '
$ws.Cells.Item(51, 5).Value = '  +4.59%  '
